# Add season record columns (Wins, Losses, Ties) to the right of the
# existing data, matching the shape/style of the existing header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy formatting from an existing header cell (bold/border/alignment)
# onto the new header cells so they share the same style index.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (same values for every player row) for
# every data row currently in the sheet.
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 63
    $ws.Cells.Item($r, 31).Value = 98
    $ws.Cells.Item($r, 32).Value = 1
}
